# Update odds/values in row 2 of the active sheet to reflect the latest
# FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 4.2
    "H2"  = 3.5
    "I2"  = 1.9
    "J2"  = 4.75
    "K2"  = 2.05
    "L2"  = 2.6
    "N2"  = 8.5
    "O2"  = 1.36
    "P2"  = 3.2
    "Q2"  = 2.2
    "R2"  = 1.67
    "S2"  = 1.44
    "T2"  = 2.63
    "W2"  = 10
    "X2"  = 21
    "Z2"  = 41
    "AC2" = 8.5
    "AD2" = 6.5
    "AG2" = 401
    "AJ2" = 9
    "AK2" = 15
    "AL2" = 17
    "AN2" = 6
    "AO2" = 23
    "AQ2" = 81
    "AV2" = 67
    "AX2" = 3.75
    "AY2" = 10
    "AZ2" = 23
    "BC2" = 201
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
